$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

$ws1.Range("F2").Value = 39
$ws1.Range("F3").Value = 21233
$ws1.Range("F5").Value = 341
$ws1.Range("F6").Value = 1127
$ws1.Range("F7").Value = 27
$ws1.Range("F8").Value = 7916
$ws1.Range("F9").Value = 551
$ws1.Range("F10").Value = 42
$ws1.Range("F12").Value = 312
$ws1.Range("F13").Value = 62
$ws1.Range("F14").Value = 187
$ws1.Range("F15").Value = 166
$ws1.Range("F18").Value = 227
$ws1.Range("F19").Value = 1358
$ws1.Range("F20").Value = 528
$ws1.Range("F21").Value = 84
$ws1.Range("F23").Value = 55
$ws1.Range("F25").Value = 83
$ws1.Range("F26").Value = 349
$ws1.Range("F27").Value = 1181
$ws1.Range("F30").Value = 222
$ws1.Range("F35").Value = 5046
$ws1.Range("F38").Value = 44
$ws1.Range("F40").Value = 13097
$ws1.Range("F41").Value = 1364
$ws1.Range("F42").Value = 134
$ws1.Range("F44").Value = 71
$ws1.Range("F45").Value = 306
$ws1.Range("F46").Value = 432
$ws1.Range("F49").Value = 102
$ws4.Range("F2").Value = 39
$ws4.Range("F3").Value = 21233
$ws4.Range("F5").Value = 1127
$ws4.Range("F6").Value = 27
$ws4.Range("F7").Value = 7916
$ws4.Range("F8").Value = 551
$ws4.Range("F9").Value = 42
$ws4.Range("F11").Value = 312
$ws4.Range("F12").Value = 62
$ws4.Range("F13").Value = 187
$ws4.Range("F14").Value = 166
$ws4.Range("F16").Value = 227
$ws4.Range("F17").Value = 1358
$ws4.Range("F18").Value = 528
$ws4.Range("F19").Value = 84
$ws4.Range("F21").Value = 55
$ws4.Range("F23").Value = 83
$ws4.Range("F24").Value = 349
$ws4.Range("F25").Value = 1181
$ws4.Range("F28").Value = 222
$ws4.Range("F35").Value = 5046
$ws4.Range("F38").Value = 44
$ws4.Range("F40").Value = 13097
$ws4.Range("F41").Value = 1364
$ws4.Range("F42").Value = 134
$ws4.Range("F44").Value = 71
$ws4.Range("F45").Value = 306
$ws4.Range("F46").Value = 432
$ws4.Range("F49").Value = 102
